$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet "LOBSTAHS_acylRanges_trimmed" right after
#    the first sheet ("LOBSTAHS_acylRanges") and before "Notes".
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$wsNotes = $wb.Worksheets.Item("Notes")
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "LOBSTAHS_acylRanges_trimmed"

# ---------------------------------------------------------------------
# 2. Header row 1: merged B1:K1, centered, same text as sheet1's B1
# ---------------------------------------------------------------------
$headerRng = $newSheet.Range("B1:K1")
$headerRng.HorizontalAlignment = -4108  # xlCenter
$newSheet.Range("B1").Value = "Total no. of fatty acid double bonds by lipid class"
$headerRng.Merge()

# ---------------------------------------------------------------------
# 3. Header row 2: copy straight from sheet1 (identical column headers)
# ---------------------------------------------------------------------
$ws1.Range("A2:K2").Copy()
$newSheet.Range("A2").PasteSpecial(-4163)  # xlPasteValuesAndNumberFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Data rows 3-66
# ---------------------------------------------------------------------
$newSheet.Cells.Item(3,1).Value = 6
$newSheet.Cells.Item(3,8).Value = 0
$newSheet.Cells.Item(3,9).Value = 2
$newSheet.Cells.Item(4,1).Value = 7
$newSheet.Cells.Item(4,8).Value = 0
$newSheet.Cells.Item(4,9).Value = 3
$newSheet.Cells.Item(5,1).Value = 8
$newSheet.Cells.Item(5,8).Value = 0
$newSheet.Cells.Item(5,9).Value = 3
$newSheet.Cells.Item(6,1).Value = 9
$newSheet.Cells.Item(6,8).Value = 0
$newSheet.Cells.Item(6,9).Value = 3
$newSheet.Cells.Item(7,1).Value = 10
$newSheet.Cells.Item(7,8).Value = 0
$newSheet.Cells.Item(7,9).Value = 4
$newSheet.Cells.Item(8,1).Value = 11
$newSheet.Cells.Item(8,8).Value = 0
$newSheet.Cells.Item(8,9).Value = 4
$newSheet.Cells.Item(9,1).Value = 12
$newSheet.Cells.Item(9,8).Value = 0
$newSheet.Cells.Item(9,9).Value = 5
$newSheet.Cells.Item(10,1).Value = 13
$newSheet.Cells.Item(11,1).Value = 14
$newSheet.Cells.Item(11,4).Value = 0
$newSheet.Cells.Item(11,5).Value = 2
$newSheet.Cells.Item(11,6).Value = 0
$newSheet.Cells.Item(11,7).Value = 2
$newSheet.Cells.Item(12,1).Value = 15
$newSheet.Cells.Item(12,4).Value = 0
$newSheet.Cells.Item(12,5).Value = 2
$newSheet.Cells.Item(12,6).Value = 0
$newSheet.Cells.Item(12,7).Value = 2
$newSheet.Cells.Item(13,1).Value = 16
$newSheet.Cells.Item(13,4).Value = 0
$newSheet.Cells.Item(13,5).Value = 4
$newSheet.Cells.Item(13,6).Value = 0
$newSheet.Cells.Item(13,7).Value = 4
$newSheet.Cells.Item(14,1).Value = 17
$newSheet.Cells.Item(14,4).Value = 0
$newSheet.Cells.Item(14,5).Value = 2
$newSheet.Cells.Item(14,6).Value = 0
$newSheet.Cells.Item(14,7).Value = 2
$newSheet.Cells.Item(15,1).Value = 18
$newSheet.Cells.Item(15,4).Value = 0
$newSheet.Cells.Item(15,5).Value = 5
$newSheet.Cells.Item(15,6).Value = 0
$newSheet.Cells.Item(15,7).Value = 5
$newSheet.Cells.Item(16,1).Value = 19
$newSheet.Cells.Item(16,4).Value = 0
$newSheet.Cells.Item(16,5).Value = 2
$newSheet.Cells.Item(16,6).Value = 0
$newSheet.Cells.Item(16,7).Value = 2
$newSheet.Cells.Item(17,1).Value = 20
$newSheet.Cells.Item(17,4).Value = 0
$newSheet.Cells.Item(17,5).Value = 5
$newSheet.Cells.Item(17,6).Value = 0
$newSheet.Cells.Item(17,7).Value = 5
$newSheet.Cells.Item(18,1).Value = 21
$newSheet.Cells.Item(18,4).Value = 0
$newSheet.Cells.Item(18,5).Value = 2
$newSheet.Cells.Item(18,6).Value = 0
$newSheet.Cells.Item(18,7).Value = 2
$newSheet.Cells.Item(19,1).Value = 22
$newSheet.Cells.Item(19,4).Value = 0
$newSheet.Cells.Item(19,5).Value = 6
$newSheet.Cells.Item(19,6).Value = 0
$newSheet.Cells.Item(19,7).Value = 6
$newSheet.Cells.Item(20,1).Value = 23
$newSheet.Cells.Item(21,1).Value = 24
$newSheet.Cells.Item(22,1).Value = 25
$newSheet.Cells.Item(23,1).Value = 26
$newSheet.Cells.Item(24,1).Value = 27
$newSheet.Cells.Item(25,1).Value = 28
$newSheet.Cells.Item(25,2).Value = 0
$newSheet.Cells.Item(25,3).Value = 7
$newSheet.Cells.Item(26,1).Value = 29
$newSheet.Cells.Item(26,2).Value = 0
$newSheet.Cells.Item(26,3).Value = 7
$newSheet.Cells.Item(27,1).Value = 30
$newSheet.Cells.Item(27,2).Value = 0
$newSheet.Cells.Item(27,3).Value = 7
$newSheet.Cells.Item(28,1).Value = 31
$newSheet.Cells.Item(28,2).Value = 0
$newSheet.Cells.Item(28,3).Value = 7
$newSheet.Cells.Item(29,1).Value = 32
$newSheet.Cells.Item(29,2).Value = 0
$newSheet.Cells.Item(29,3).Value = 8
$newSheet.Cells.Item(30,1).Value = 33
$newSheet.Cells.Item(30,2).Value = 0
$newSheet.Cells.Item(30,3).Value = 8
$newSheet.Cells.Item(31,1).Value = 34
$newSheet.Cells.Item(31,2).Value = 0
$newSheet.Cells.Item(31,3).Value = 9
$newSheet.Cells.Item(32,1).Value = 35
$newSheet.Cells.Item(32,2).Value = 0
$newSheet.Cells.Item(32,3).Value = 8
$newSheet.Cells.Item(33,1).Value = 36
$newSheet.Cells.Item(33,2).Value = 0
$newSheet.Cells.Item(33,3).Value = 10
$newSheet.Cells.Item(34,1).Value = 37
$newSheet.Cells.Item(34,2).Value = 0
$newSheet.Cells.Item(34,3).Value = 8
$newSheet.Cells.Item(35,1).Value = 38
$newSheet.Cells.Item(35,2).Value = 0
$newSheet.Cells.Item(35,3).Value = 10
$newSheet.Cells.Item(36,1).Value = 39
$newSheet.Cells.Item(36,2).Value = 0
$newSheet.Cells.Item(36,3).Value = 8
$newSheet.Cells.Item(37,1).Value = 40
$newSheet.Cells.Item(37,2).Value = 0
$newSheet.Cells.Item(37,3).Value = 11
$newSheet.Cells.Item(38,1).Value = 41
$newSheet.Cells.Item(38,2).Value = 0
$newSheet.Cells.Item(38,3).Value = 8
$newSheet.Cells.Item(39,1).Value = 42
$newSheet.Cells.Item(39,2).Value = 0
$newSheet.Cells.Item(39,3).Value = 11
$newSheet.Cells.Item(39,10).Value = 0
$newSheet.Cells.Item(39,11).Value = 10
$newSheet.Cells.Item(40,1).Value = 43
$newSheet.Cells.Item(40,2).Value = 0
$newSheet.Cells.Item(40,3).Value = 8
$newSheet.Cells.Item(40,10).Value = 0
$newSheet.Cells.Item(40,11).Value = 10
$newSheet.Cells.Item(41,1).Value = 44
$newSheet.Cells.Item(41,2).Value = 0
$newSheet.Cells.Item(41,3).Value = 12
$newSheet.Cells.Item(41,10).Value = 0
$newSheet.Cells.Item(41,11).Value = 11
$newSheet.Cells.Item(42,1).Value = 45
$newSheet.Cells.Item(42,10).Value = 0
$newSheet.Cells.Item(42,11).Value = 11
$newSheet.Cells.Item(43,1).Value = 46
$newSheet.Cells.Item(43,10).Value = 0
$newSheet.Cells.Item(43,11).Value = 12
$newSheet.Cells.Item(44,1).Value = 47
$newSheet.Cells.Item(44,10).Value = 0
$newSheet.Cells.Item(44,11).Value = 12
$newSheet.Cells.Item(45,1).Value = 48
$newSheet.Cells.Item(45,10).Value = 0
$newSheet.Cells.Item(45,11).Value = 12
$newSheet.Cells.Item(46,1).Value = 49
$newSheet.Cells.Item(46,10).Value = 0
$newSheet.Cells.Item(46,11).Value = 12
$newSheet.Cells.Item(47,1).Value = 50
$newSheet.Cells.Item(47,10).Value = 0
$newSheet.Cells.Item(47,11).Value = 13
$newSheet.Cells.Item(48,1).Value = 51
$newSheet.Cells.Item(48,10).Value = 0
$newSheet.Cells.Item(48,11).Value = 13
$newSheet.Cells.Item(49,1).Value = 52
$newSheet.Cells.Item(49,10).Value = 0
$newSheet.Cells.Item(49,11).Value = 14
$newSheet.Cells.Item(50,1).Value = 53
$newSheet.Cells.Item(50,10).Value = 0
$newSheet.Cells.Item(50,11).Value = 14
$newSheet.Cells.Item(51,1).Value = 54
$newSheet.Cells.Item(51,10).Value = 0
$newSheet.Cells.Item(51,11).Value = 15
$newSheet.Cells.Item(52,1).Value = 55
$newSheet.Cells.Item(52,10).Value = 0
$newSheet.Cells.Item(52,11).Value = 14
$newSheet.Cells.Item(53,1).Value = 56
$newSheet.Cells.Item(53,10).Value = 0
$newSheet.Cells.Item(53,11).Value = 15
$newSheet.Cells.Item(54,1).Value = 57
$newSheet.Cells.Item(54,10).Value = 0
$newSheet.Cells.Item(54,11).Value = 14
$newSheet.Cells.Item(55,1).Value = 58
$newSheet.Cells.Item(55,10).Value = 0
$newSheet.Cells.Item(55,11).Value = 8
$newSheet.Cells.Item(56,1).Value = 59
$newSheet.Cells.Item(57,1).Value = 60
$newSheet.Cells.Item(58,1).Value = 61
$newSheet.Cells.Item(59,1).Value = 62
$newSheet.Cells.Item(60,1).Value = 63
$newSheet.Cells.Item(61,1).Value = 64
$newSheet.Cells.Item(62,1).Value = 65
$newSheet.Cells.Item(63,1).Value = 66
$newSheet.Cells.Item(64,1).Value = 67
$newSheet.Cells.Item(65,1).Value = 68
$newSheet.Cells.Item(66,1).Value = 69

# ---------------------------------------------------------------------
# 5. Notes sheet: add row 18 documenting this change
# ---------------------------------------------------------------------
$wsNotes.Range("A18").Value = 43025
$ws1.Range("A14").Copy()
$wsNotes.Range("A18").PasteSpecial(-4122)  # xlPasteFormats (reuse date style)
$excel.CutCopyMode = 0
$wsNotes.Range("B18").Value = "Added an alternate acylRanges table (""LOBSTAHS_acylRanges_trimmed"") that will generate a reduced database of species with a narrower range of acyl chain lengths; for use in the Van Mooy Lab"
$wsNotes.Range("C18").Value = "TL"

# ---------------------------------------------------------------------
# 6. View-state bits present in the diff
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B1:K1").Select()

$newSheet.Activate()
$newSheet.Range("G18").Select()

$wsNotes.Activate()
$wsNotes.Range("A19").Select()

Write-Output "done"
